$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C5").Value = 4
$ws.Range("D5").Value = "2026-02-12T11:43:30.716466+00:00"
$ws.Range("E5").Value = "Ohkkkkk Omg thanks bro 😭😭 It means alot Btw do u know any software for cheating in oa? How? Tell tell Hackerank pe bhi? Session of what? Omg I'm confused Is it free or paid? How much He who???? Isn't that ai based? Ohhhh But what if I give it to my frnd? ?? Ohh u are talking about remote access thing?? For lifetime? Oh Ok"
$ws.Range("G5").Value = "Dude, chill, it's not for cheating, session means proctoring, and yeah, it's paid, don't share with friends btw 😅"
